# Update countries & provincias Spain
#
# This script applies the 16-Apr-2020 data refresh to the "paises" sheet:
#  - Updates case/death figures for a handful of countries whose row
#    position does not change (Estados Unidos, Alemania, Finlandia,
#    Argentina, Vietnam).
#  - Re-orders three pairs of countries (India before Israel, Azerbaiyan
#    before Eslovenia, Kenia before Mayotte) and refreshes the figures
#    for the country that moves into the new slot, while the countries
#    that shift down one row keep their previous figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    # $values is an array: Country, CasosTotales, NuevosCasos, CasosActivos,
    # Recuperados, CasosCriticos, MuertesHoy, Muertes
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

# ---------------------------------------------------------------------
# Simple in-place value refreshes (row/country unchanged)
# ---------------------------------------------------------------------

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 644806
$ws.Cells.Item(4, 3).Value = 717
$ws.Cells.Item(4, 4).Value = 48710
$ws.Cells.Item(4, 5).Value = 567524
$ws.Cells.Item(4, 7).Value = 43
$ws.Cells.Item(4, 8).Value = 28572

# Row 8: Alemania
$ws.Cells.Item(8, 2).Value = 135230
$ws.Cells.Item(8, 3).Value = 477
$ws.Cells.Item(8, 5).Value = 54381
$ws.Cells.Item(8, 7).Value = 45
$ws.Cells.Item(8, 8).Value = 3849

# Row 51: Finlandia
$ws.Cells.Item(51, 5).Value = 2994
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 75

# Row 54: Argentina
$ws.Cells.Item(54, 6).Value = 121

# Row 115: Vietnam
$ws.Cells.Item(115, 4).Value = 177
$ws.Cells.Item(115, 5).Value = 91

# ---------------------------------------------------------------------
# Re-order block 1: India moves to just above Israel (rows 21-24)
#   before: Israel(21) Irlanda(22) Suecia(23) India(24)
#   after : India(21)  Israel(22)  Irlanda(23) Suecia(24)
# ---------------------------------------------------------------------

Set-Row 21 @("India",   12759, 389, 1514,  10822, 0,   1,   423)
Set-Row 22 @("Israel",  12591, 90,  2624,  9827,  174, 10,  140)
Set-Row 23 @("Irlanda", 12547, 0,   77,    12026, 158, 0,   444)
Set-Row 24 @("Suecia",  12540, 613, 381,   10826, 996, 130, 1333)

# ---------------------------------------------------------------------
# Re-order block 2: Azerbaiyan moves to just above Eslovenia (rows 72-73)
#   before: Eslovenia(72) Azerbaiyan(73)
#   after : Azerbaiyan(72) Eslovenia(73)
# ---------------------------------------------------------------------

Set-Row 72 @("Azerbaiyan", 1283, 30, 460, 808,  28, 2, 15)
Set-Row 73 @("Eslovenia",  1268, 20, 174, 1033, 31, 0, 61)

# ---------------------------------------------------------------------
# Re-order block 3: Kenia moves to just above Mayotte (rows 118-119)
#   before: Mayotte(118) Kenia(119)
#   after : Kenia(118)   Mayotte(119)
# ---------------------------------------------------------------------

Set-Row 118 @("Kenia",   234, 9,  53, 170, 2, 1, 11)
Set-Row 119 @("Mayotte", 233, 16, 69, 161, 3, 0, 3)

$wb.Save()
